# edit.ps1 - applies the "ExcelExport korrigiert" fix:
#  1. Renames column header F1 from "DauerStunden" to "DauerMinuten"
#     (the column always held minutes, the old header text was wrong)
#  2. Appends the newly scanned arrival/departure log rows (30-46)
#     that were exported since the last run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix mislabeled header -------------------------------------------
$ws.Cells.Item(1, 6).Value = "DauerMinuten"

# --- 2) Append new log rows ----------------------------------------------
$newRows = @(
    ,@("Eli", "Enders", "4a", "13.11.2023 21:49", "13.11.2023 21:49", 0)
    ,@("Max", "Schmitz", "4a", "13.11.2023 21:51", "13.11.2023 21:51", 0)
    ,@("Stephan", "Fuchs", "3C", "13.11.2023 21:58", "13.11.2023 21:58", 0)
    ,@("Stephan", "Fuchs", "3C", "14.11.2023 17:29", "14.11.2023 17:29", 0)
    ,@("Max", "Schmitz", "4a", "19.11.2023 18:54", "19.11.2023 19:34", 40)
    ,@("Eli", "Enders", "4a", "19.11.2023 19:35", "19.11.2023 19:35", 0)
    ,@("Eli", "Enders", "4a", "22.11.2023 10:32", "22.11.2023 10:32", 0)
    ,@("Detlef", "Soost", "1a", "22.11.2023 13:53", "22.11.2023 13:53", 0)
    ,@("Stephan", "Fuchs", "3C", "26.11.2023 16:02", "26.11.2023 16:02", 0)
    ,@("Stephan", "Fuchs", "3C", "07.12.2023 20:40", "07.12.2023 20:41", 1)
    ,@("Eli", "Enders", "4a", "07.12.2023 20:40", "07.12.2023 20:41", 1)
    ,@("Stephan", "Fuchs", "3C", "17.12.2023 22:06", "17.12.2023 22:06", 0)
    ,@("Eli", "Enders", "4a", "17.12.2023 22:06", "17.12.2023 22:06", 0)
    ,@("Detlef", "Soost", "1a", "17.12.2023 22:06", "17.12.2023 22:06", 0)
    ,@("Stephan", "Fuchs", "3C", "18.12.2023 08:00", "18.12.2023 08:08", 8)
    ,@("Max", "Schmitz", "4a", "18.12.2023 08:00", "18.12.2023 08:01", 1)
    ,@("Detlef", "Soost", "1a", "18.12.2023 08:00", "18.12.2023 08:01", 1)
)

$startRow = 30
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
